$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = "diya"
$ws.Range("C3").Value = "cs22B"
$ws.Range("D3").Value = "86h"
$ws.Range("E3").Value = "diya"
$ws.Range("F3").Value = (Get-Date -Year 2019 -Month 5 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G3").Value = (Get-Date -Year 2019 -Month 5 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("H3").Value = (Get-Date -Year 2019 -Month 5 -Day 22 -Hour 0 -Minute 0 -Second 0)
